# Updated cryptos list on Sat Oct 26 19:52:57 UTC 2024 with GitHub Actions
#
# D/E columns hold plain text (inline strings), even when the text looks
# like a plain decimal number (e.g. "585.07"). Writing such a value through
# Range.Value would make Excel auto-detect it as a numeric literal, so for
# those specific cells we force a text literal (NumberFormat "@", write the
# value, then restore the "Normal" style so no stray per-cell formatting is
# left behind) and for everything else (values that aren't plain numbers,
# like "67.074.47" or "  +0.15%  ") a normal assignment already keeps them
# as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

function Set-Cell {
    param($row, $col, [string]$value)
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 2 - Bitcoin
Set-Cell 2 4 "67.074.47"
Set-Cell 2 5 "  +0.15%  "

# Row 3 - Ethereum
Set-Cell 3 4 "2.482.05"
Set-Cell 3 5 "  +0.20%  "

# Row 4 - TetherUSD
Set-Cell 4 5 "  -0.02%  "

# Row 5 - BNB
Set-TextValue $ws.Cells.Item(5, 4) "585.07"
Set-Cell 5 5 "  -0.03%  "

# Row 6 - Solana
Set-TextValue $ws.Cells.Item(6, 4) "171.26"
Set-Cell 6 5 "  +1.86%  "

# Row 8 - XRP
Set-Cell 8 5 "  -0.31%  "

# Row 9 - LidoStakedEther
Set-Cell 9 4 "2.481.92"
Set-Cell 9 5 "  +0.07%  "

# Row 10 - Dogecoin
Set-Cell 10 5 "  +0.80%  "

# Row 11 - TRON
Set-Cell 11 5 "  +0.03%  "

# Row 12 - Toncoin
Set-Cell 12 5 "  -0.35%  "

# Row 13 - Cardano
Set-Cell 13 5 "  -1.41%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-Cell 14 5 "  +0.04%  "

# Row 15 - Avalanche
Set-TextValue $ws.Cells.Item(15, 4) "25.42"
Set-Cell 15 5 "  -1.84%  "

# Row 16 - WrappedBTC
Set-Cell 16 4 "66.955.92"
Set-Cell 16 5 "  +0.14%  "

# Row 17 - ShibaInu
Set-Cell 17 5 "  -1.80%  "

# Row 18 - WrappedEther
Set-Cell 18 4 "2.480.18"
Set-Cell 18 5 "  -0.46%  "

# Row 19 - Chainlink
Set-Cell 19 5 "  -5.32%  "

# Row 20 - Uniswap
Set-Cell 20 5 "  -5.91%  "

# Row 21 - BitcoinCash
Set-Cell 21 5 "  -3.39%  "

# Row 23 - Dai
Set-Cell 23 5 "  +0.19%  "

# Row 24 - Litecoin
Set-TextValue $ws.Cells.Item(24, 4) "68.48"
Set-Cell 24 5 "  -3.31%  "

# Row 25 - NEARProtocol
Set-Cell 25 5 "  -4.86%  "

# Row 26 - SuiNetwork
Set-Cell 26 5 "  -2.35%  "

# Row 27 - Aptos
Set-Cell 27 5 "  -1.96%  "

# Row 28 - Binance-PegBSC-USD
Set-TextValue $ws.Cells.Item(28, 4) "0.998"
Set-Cell 28 5 "  -0.41%  "

# Row 29 - WrappedeETH
Set-Cell 29 4 "2.607.17"
Set-Cell 29 5 "  -0.45%  "

# Row 30 - PEPE
Set-Cell 30 4 "0.0₃0899"
Set-Cell 30 5 "  -2.91%  "

# Row 31 - Bittensor
Set-TextValue $ws.Cells.Item(31, 4) "510.35"
Set-Cell 31 5 "  -0.86%  "

# Row 32 - InternetComputer(DFINITY)
Set-Cell 32 5 "  -5.09%  "

# Row 33 - Fetch.AI
Set-Cell 33 5 "  -3.16%  "

# Row 34 - PancakeSwap
Set-Cell 34 5 "  -3.79%  "

# Row 35 - FirstDigitalUSD
Set-Cell 35 5 "  -0.05%  "

# Row 36 - Monero
Set-TextValue $ws.Cells.Item(36, 4) "159.86"
Set-Cell 36 5 "  +1.97%  "

# Row 38 - WhiteBITCoin
Set-Cell 38 5 "  +0.65%  "

# Row 39 - EthereumClassic
Set-Cell 39 5 "  -3.63%  "

# Row 40 - ImmutableX
Set-Cell 40 5 "  -5.93%  "

# Row 42 - Stacks
Set-Cell 42 5 "  -3.46%  "

# Row 43 / Row 44 - RenderToken and PolygonEcosystemToken swap position
Set-Cell 43 2 "PolygonEcosystemToken"
Set-Cell 43 3 "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
Set-TextValue $ws.Cells.Item(43, 4) "0.327"
Set-Cell 43 5 "  -1.56%  "

Set-Cell 44 2 "RenderToken"
Set-Cell 44 3 "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue $ws.Cells.Item(44, 4) "4.80"
Set-Cell 44 5 "  -2.88%  "

# Row 45 - dogwifhat
Set-Cell 45 5 "  -3.14%  "

# Row 46 - OKB
Set-TextValue $ws.Cells.Item(46, 4) "38.77"
Set-Cell 46 5 "  -1.06%  "

# Row 47 - Aave
Set-TextValue $ws.Cells.Item(47, 4) "142.73"
Set-Cell 47 5 "  +0.01%  "

# Row 48 - ARBITRUM
Set-Cell 48 5 "  -4.16%  "

# Row 49 - Filecoin
Set-Cell 49 5 "  -4.30%  "

# Row 50 - BabyDogeCoin
Set-Cell 50 5 "  -6.35%  "

# Row 51 - Cronos
Set-Cell 51 5 "  -0.93%  "
